$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.826.18"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "3.318.35"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'574.59"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").Value = "'182.57"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.602"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "3.316.66"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "'0.129"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "'0.403"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "3.897.67"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "'27.08"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "67.004.71"
$ws.Range("E16").Value = "  -1.74%  "
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "3.324.39"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").Value = "'439.79"
$ws.Range("E19").Value = "  +5.59%  "
$ws.Range("D20").Value = "'13.53"
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").Value = "'73.85"
$ws.Range("E23").Value = "  +3.86%  "
$ws.Range("D24").Value = "'0.996"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").Value = "3.474.10"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("D26").Value = "'0.511"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E28").Value = "  +3.71%  "
$ws.Range("D29").Value = "'8.92"
$ws.Range("E29").Value = "  -4.64%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").Value = "'22.88"
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "'5.28"
$ws.Range("E34").Value = "  -3.15%  "
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("D36").Value = "'1.21"
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("E37").Value = "  +3.75%  "
$ws.Range("D38").Value = "'161.22"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").Value = "'27.47"
$ws.Range("E39").Value = "  +3.11%  "
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("D41").Value = "2.814.81"
$ws.Range("E41").Value = "  +6.67%  "
$ws.Range("D42").Value = "'0.787"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "'4.44"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").Value = "'40.32"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "'2.28"
$ws.Range("E48").Value = "  -5.27%  "
$ws.Range("D49").Value = "'319.09"
$ws.Range("E49").Value = "  -4.65%  "
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("E51").Value = "  -1.16%  "
